$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999998819767721
$ws.Range("A2").Value = 0.9936527946818593
$ws.Range("A3").Value = 0.97240432223391615
$ws.Range("A4").Value = 0.96259902338472525
$ws.Range("A5").Value = 0.95320305584342291
$ws.Range("A6").Value = 0.93041519461904976
$ws.Range("A7").Value = 0.92643697945687586
$ws.Range("A8").Value = 0.92196542484724231
$ws.Range("A9").Value = 0.91788593917980266
$ws.Range("A10").Value = 0.91473998218573793
$ws.Range("A11").Value = 0.91440647202058911
$ws.Range("A12").Value = 0.91402783403889387
$ws.Range("A13").Value = 0.9137086751626633
$ws.Range("A14").Value = 0.91424553197912117
$ws.Range("A15").Value = 0.91618723104919608
$ws.Range("A16").Value = 0.91912417578916061
$ws.Range("A17").Value = 0.9266663702584792
$ws.Range("A18").Value = 0.92695029278508445
$ws.Range("A19").Value = 0.98891379108809607
$ws.Range("A20").Value = 0.96447946096236326
$ws.Range("A21").Value = 0.9578581776809858
$ws.Range("A22").Value = 0.95659367065923739
$ws.Range("A23").Value = 0.96536878336650012
$ws.Range("A24").Value = 0.95234713220893608
$ws.Range("A25").Value = 0.94589002578555204
$ws.Range("A26").Value = 0.92786838632990087
$ws.Range("A27").Value = 0.92456677650037777
$ws.Range("A28").Value = 0.91047282748616465
$ws.Range("A29").Value = 0.89520171603043464
$ws.Range("A30").Value = 0.88863146616076416
$ws.Range("A31").Value = 0.88097801161101097
$ws.Range("A32").Value = 0.87929876697706311
$ws.Range("A33").Value = 0.87877878169975832
